$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.802.09"
$ws.Range("E2").Value = "  -2.69%  "

# Row 3
$ws.Range("D3").Value = "3.186.08"
$ws.Range("E3").Value = "  -1.46%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "'600.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("D6").Value = "'152.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.17%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "3.184.89"
$ws.Range("E8").Value = "  -1.56%  "

# Row 9
$ws.Range("D9").Value = "'0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.63%  "

# Row 10
$ws.Range("D10").Value = "'0.153"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.59%  "

# Row 11
$ws.Range("D11").Value = "'5.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "

# Row 12
$ws.Range("D12").Value = "'0.477"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.92%  "

# Row 13
$ws.Range("E13").Value = "  -6.41%  "

# Row 14
$ws.Range("D14").Value = "'36.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.67%  "

# Row 15
$ws.Range("D15").Value = "3.696.45"
$ws.Range("E15").Value = "  -1.63%  "

# Row 16
$ws.Range("D16").Value = "64.873.47"
$ws.Range("E16").Value = "  -2.59%  "

# Row 17
$ws.Range("D17").Value = "3.182.83"
$ws.Range("E17").Value = "  -1.42%  "

# Row 18
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("E19").Value = "  -5.01%  "

# Row 20
$ws.Range("D20").Value = "'481.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.72%  "

# Row 21
$ws.Range("D21").Value = "'14.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.30%  "

# Row 23
$ws.Range("D23").Value = "'7.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.73%  "

# Row 24
$ws.Range("D24").Value = "'13.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.86%  "

# Row 25
$ws.Range("D25").Value = "'84.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").Value = "'2.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.05%  "

# Row 28
$ws.Range("D28").Value = "'8.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.75%  "

# Row 29
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.49%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.126"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +25.91%  "

# Row 31
$ws.Range("D31").Value = "'6.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.35%  "

# Row 32
$ws.Range("E32").Value = "  -8.78%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'26.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.60%  "

# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "

# Row 35
$ws.Range("D35").Value = "'1.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.39%  "

# Row 36
$ws.Range("D36").Value = "'6.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.97%  "

# Row 37
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'54.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.95%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.47%  "

# Row 39
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0738"
$ws.Range("E39").Value = "  -4.88%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'459.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.18%  "

# Row 41
$ws.Range("D41").Value = "'0.0403"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.49%  "

# Row 42
$ws.Range("E42").Value = "  -3.81%  "

# Row 43
$ws.Range("D43").Value = "'8.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.00%  "

# Row 44
$ws.Range("E44").Value = "  -2.15%  "

# Row 45
$ws.Range("D45").Value = "2.893.54"
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("D46").Value = "'0.276"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.94%  "

# Row 47
$ws.Range("D47").Value = "'27.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.67%  "

# Row 48
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.96%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.116"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("D51").Value = "'120.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
